# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# (and the Chainlink/Polygon/Polkadot/WrappedEther row re-ordering + Gas/RenderToken swap)
# described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.476.86"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "1.880.83"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Value = "'240.68"
$ws.Range("E5").Value = "  +3.94%  "
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").Value = "'42.99"
$ws.Range("E8").Value = "  +7.75%  "
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "2.149.47"
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.908.25"
$ws.Range("E13").Value = "  +3.52%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.69"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.688"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'4.76"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").Value = "35.463.55"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").Value = "'70.68"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").Value = "'242.79"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "'12.43"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("D22").Value = "'4.78"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").Value = "'170.25"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'1.92"
$ws.Range("E26").Value = "  +26.33%  "
$ws.Range("D27").Value = "'8.28"
$ws.Range("E27").Value = "  +6.22%  "
$ws.Range("D28").Value = "'17.88"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  +25.29%  "
$ws.Range("E35").Value = "  +7.95%  "
$ws.Range("D36").Value = "'0.829"
$ws.Range("E37").Value = "  +7.07%  "
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("E39").Value = "  +5.16%  "
$ws.Range("D40").Value = "'91.57"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "1.356.29"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "'15.33"
$ws.Range("E42").Value = "  +3.32%  "
$ws.Range("D43").Value = "'0.0604"
$ws.Range("E43").Value = "  +15.47%  "
$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").Value = "'13.16"
$ws.Range("E44").Value = "  +58.01%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.37"
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("D46").Value = "'2.42"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  +6.40%  "
$ws.Range("D48").Value = "'2.72"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "2.063.52"
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").Value = "'0.0688"
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("D51").Value = "'3.46"
$ws.Range("E51").Value = "  -0.58%  "
